$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 04:05"

# Row 67
$ws.Range("B67").Value = 5187
$ws.Range("C67").Value = 268
$ws.Range("D67").Value = 561
$ws.Range("E67").Value = 4411
$ws.Range("G67").Value = 16
$ws.Range("H67").Value = 215

# Row 80
$ws.Range("A80").Value = "Guatemala"
$ws.Range("B80").Value = 2512
$ws.Range("C80").Value = 247
$ws.Range("D80").Value = 222
$ws.Range("E80").Value = 2242
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 48

# Row 81
$ws.Range("A81").Value = "Tayikistan"
$ws.Range("D81").Value = 470
$ws.Range("E81").Value = 1836
$ws.Range("H81").Value = 44

# Row 82
$ws.Range("A82").Value = "Bosnia y Herzegovina"
$ws.Range("B82").Value = 2350
$ws.Range("D82").Value = 1596
$ws.Range("E82").Value = 614
$ws.Range("H82").Value = 140

# Row 83
$ws.Range("A83").Value = "Bulgaria"
$ws.Range("B83").Value = 2331
$ws.Range("D83").Value = 727
$ws.Range("E83").Value = 1484
$ws.Range("H83").Value = 120

# Row 84
$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 2301
$ws.Range("D84").Value = 1100
$ws.Range("E84").Value = 1172
$ws.Range("H84").Value = 29

# Row 121
$ws.Range("A121").Value = "Haiti"
$ws.Range("B121").Value = 734
$ws.Range("C121").Value = 71
$ws.Range("D121").Value = 21
$ws.Range("E121").Value = 688
$ws.Range("G121").Value = 3
$ws.Range("H121").Value = 25

# Row 122
$ws.Range("A122").Value = "Georgia"
$ws.Range("B122").Value = 721
$ws.Range("D122").Value = 485
$ws.Range("E122").Value = 224
$ws.Range("H122").Value = 12

# Row 123
$ws.Range("A123").Value = "Crucero"
$ws.Range("B123").Value = 712
$ws.Range("D123").Value = 651
$ws.Range("E123").Value = 48
$ws.Range("H123").Value = 13

# Row 124
$ws.Range("A124").Value = "Jordania"
$ws.Range("B124").Value = 684
$ws.Range("D124").Value = 457
$ws.Range("E124").Value = 218
$ws.Range("H124").Value = 9

# Row 144
$ws.Range("D144").Value = 266
$ws.Range("E144").Value = 58

# Row 209
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

# Row 210
$ws.Range("A210").Value = "Groenlandia"

# Row 211
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 214
$ws.Range("A214").Value = "Sahara Occidental"

# Row 215
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
